# Update NATMI ligand-receptor edge table (Col2a1-Itga2) with values
# recomputed from the new TPM expression matrix.
#
# The workbook structure is unchanged; only the per-row ligand/receptor
# expression, specificity and edge-weight statistics (columns E-T) are
# refreshed to reflect the regenerated TPM data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 0.05522166666666666
$ws.Range("H2").Value = 0.165665
$ws.Range("I2").Value = 0.1626033416795164
$ws.Range("J2").Value = 0.1626033416795164
$ws.Range("M2").Value = 3.339352
$ws.Range("N2").Value = 10.018056
$ws.Range("O2").Value = 0.6054960700393903
$ws.Range("P2").Value = 0.6054960700393903
$ws.Range("Q2").Value = 0.1844045830266667
$ws.Range("R2").Value = 1.65964124724
$ws.Range("S2").Value = 0.09845568436221937
$ws.Range("T2").Value = 0.09845568436221937

# Row 3
$ws.Range("G3").Value = 0.05522166666666666
$ws.Range("H3").Value = 0.165665
$ws.Range("I3").Value = 0.1626033416795164
$ws.Range("J3").Value = 0.1626033416795164
$ws.Range("O3").Value = 0.2540955070726236
$ws.Range("P3").Value = 0.2540955070726236
$ws.Range("Q3").Value = 0.07738510347000001
$ws.Range("R3").Value = 0.69646593123
$ws.Range("S3").Value = 0.04131677855575979
$ws.Range("T3").Value = 0.04131677855575978

# Row 4
$ws.Range("G4").Value = 0.05522166666666666
$ws.Range("H4").Value = 0.165665
$ws.Range("I4").Value = 0.1626033416795164
$ws.Range("J4").Value = 0.1626033416795164
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1338136666666667
$ws.Range("N4").Value = 0.401441
$ws.Range("O4").Value = 0.02426328499787613
$ws.Range("P4").Value = 0.02426328499787612
$ws.Range("Q4").Value = 0.007389413696111112
$ws.Range("R4").Value = 0.066504723265
$ws.Range("S4").Value = 0.003945291220777135
$ws.Range("T4").Value = 0.003945291220777134

# Row 5
$ws.Range("G5").Value = 0.05522166666666666
$ws.Range("H5").Value = 0.165665
$ws.Range("I5").Value = 0.1626033416795164
$ws.Range("J5").Value = 0.1626033416795164
$ws.Range("M5").Value = 0.6405483333333334
$ws.Range("N5").Value = 1.921645
$ws.Range("O5").Value = 0.11614513789011
$ws.Range("P5").Value = 0.11614513789011
$ws.Range("Q5").Value = 0.03537214654722222
$ws.Range("R5").Value = 0.318349318925
$ws.Range("S5").Value = 0.01888558754076011
$ws.Range("T5").Value = 0.0188855875407601

# Row 6
$ws.Range("I6").Value = 0.4895738146440669
$ws.Range("J6").Value = 0.4895738146440669
$ws.Range("M6").Value = 3.339352
$ws.Range("N6").Value = 10.018056
$ws.Range("O6").Value = 0.6054960700393903
$ws.Range("P6").Value = 0.6054960700393903
$ws.Range("Q6").Value = 0.555214020928
$ws.Range("R6").Value = 4.996926188352001
$ws.Range("S6").Value = 0.2964350207611754
$ws.Range("T6").Value = 0.2964350207611754

# Row 7
$ws.Range("I7").Value = 0.4895738146440669
$ws.Range("J7").Value = 0.4895738146440669
$ws.Range("O7").Value = 0.2540955070726236
$ws.Range("P7").Value = 0.2540955070726236
$ws.Range("S7").Value = 0.1243985066814628
$ws.Range("T7").Value = 0.1243985066814628

# Row 8
$ws.Range("I8").Value = 0.4895738146440669
$ws.Range("J8").Value = 0.4895738146440669
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.1338136666666667
$ws.Range("N8").Value = 0.401441
$ws.Range("O8").Value = 0.02426328499787613
$ws.Range("P8").Value = 0.02426328499787612
$ws.Range("Q8").Value = 0.02224839547466667
$ws.Range("R8").Value = 0.200235559272
$ws.Range("S8").Value = 0.01187866899220637
$ws.Range("T8").Value = 0.01187866899220637

# Row 9
$ws.Range("I9").Value = 0.4895738146440669
$ws.Range("J9").Value = 0.4895738146440669
$ws.Range("M9").Value = 0.6405483333333334
$ws.Range("N9").Value = 1.921645
$ws.Range("O9").Value = 0.11614513789011
$ws.Range("P9").Value = 0.11614513789011
$ws.Range("Q9").Value = 0.1065001280933333
$ws.Range("R9").Value = 0.95850115284
$ws.Range("S9").Value = 0.05686161820922232
$ws.Range("T9").Value = 0.05686161820922232

# Row 10
$ws.Range("G10").Value = 0.02740366666666667
$ws.Range("H10").Value = 0.08221100000000001
$ws.Range("I10").Value = 0.08069165679422162
$ws.Range("J10").Value = 0.08069165679422161
$ws.Range("M10").Value = 3.339352
$ws.Range("N10").Value = 10.018056
$ws.Range("O10").Value = 0.6054960700393903
$ws.Range("P10").Value = 0.6054960700393903
$ws.Range("Q10").Value = 0.09151048909066668
$ws.Range("R10").Value = 0.8235944018160002
$ws.Range("S10").Value = 0.04885848107386846
$ws.Range("T10").Value = 0.04885848107386845

# Row 11
$ws.Range("G11").Value = 0.02740366666666667
$ws.Range("H11").Value = 0.08221100000000001
$ws.Range("I11").Value = 0.08069165679422162
$ws.Range("J11").Value = 0.08069165679422161
$ws.Range("O11").Value = 0.2540955070726236
$ws.Range("P11").Value = 0.2540955070726236
$ws.Range("Q11").Value = 0.03840223789800001
$ws.Range("R11").Value = 0.3456201410820001
$ws.Range("S11").Value = 0.02050338744965786
$ws.Range("T11").Value = 0.02050338744965785

# Row 12
$ws.Range("G12").Value = 0.02740366666666667
$ws.Range("H12").Value = 0.08221100000000001
$ws.Range("I12").Value = 0.08069165679422162
$ws.Range("J12").Value = 0.08069165679422161
$ws.Range("K12").Value = 2
$ws.Range("L12").Value = 0.6666666666666666
$ws.Range("M12").Value = 0.1338136666666667
$ws.Range("N12").Value = 0.401441
$ws.Range("O12").Value = 0.02426328499787613
$ws.Range("P12").Value = 0.02426328499787612
$ws.Range("Q12").Value = 0.003666985116777779
$ws.Range("R12").Value = 0.03300286605100001
$ws.Range("S12").Value = 0.001957844665749007
$ws.Range("T12").Value = 0.001957844665749006

# Row 13
$ws.Range("G13").Value = 0.02740366666666667
$ws.Range("H13").Value = 0.08221100000000001
$ws.Range("I13").Value = 0.08069165679422162
$ws.Range("J13").Value = 0.08069165679422161
$ws.Range("M13").Value = 0.6405483333333334
$ws.Range("N13").Value = 1.921645
$ws.Range("O13").Value = 0.11614513789011
$ws.Range("P13").Value = 0.11614513789011
$ws.Range("Q13").Value = 0.01755337301055556
$ws.Range("R13").Value = 0.157980357095
$ws.Range("S13").Value = 0.009371943604946304
$ws.Range("T13").Value = 0.009371943604946302

# Row 14
$ws.Range("G14").Value = 0.046975
$ws.Range("H14").Value = 0.140925
$ws.Range("I14").Value = 0.1383205621355497
$ws.Range("J14").Value = 0.1383205621355497
$ws.Range("M14").Value = 3.339352
$ws.Range("N14").Value = 10.018056
$ws.Range("O14").Value = 0.6054960700393903
$ws.Range("P14").Value = 0.6054960700393903
$ws.Range("Q14").Value = 0.1568660602
$ws.Range("R14").Value = 1.4117945418
$ws.Range("S14").Value = 0.08375255677871465
$ws.Range("T14").Value = 0.08375255677871465

# Row 15
$ws.Range("G15").Value = 0.046975
$ws.Range("H15").Value = 0.140925
$ws.Range("I15").Value = 0.1383205621355497
$ws.Range("J15").Value = 0.1383205621355497
$ws.Range("O15").Value = 0.2540955070726236
$ws.Range("P15").Value = 0.2540955070726236
$ws.Range("Q15").Value = 0.06582860415
$ws.Range("R15").Value = 0.59245743735
$ws.Range("S15").Value = 0.03514663337440285
$ws.Range("T15").Value = 0.03514663337440284

# Row 16
$ws.Range("G16").Value = 0.046975
$ws.Range("H16").Value = 0.140925
$ws.Range("I16").Value = 0.1383205621355497
$ws.Range("J16").Value = 0.1383205621355497
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.1338136666666667
$ws.Range("N16").Value = 0.401441
$ws.Range("O16").Value = 0.02426328499787613
$ws.Range("P16").Value = 0.02426328499787612
$ws.Range("Q16").Value = 0.006285896991666667
$ws.Range("R16").Value = 0.056573072925
$ws.Range("S16").Value = 0.003356111220161276
$ws.Range("T16").Value = 0.003356111220161275

# Row 17
$ws.Range("G17").Value = 0.046975
$ws.Range("H17").Value = 0.140925
$ws.Range("I17").Value = 0.1383205621355497
$ws.Range("J17").Value = 0.1383205621355497
$ws.Range("M17").Value = 0.6405483333333334
$ws.Range("N17").Value = 1.921645
$ws.Range("O17").Value = 0.11614513789011
$ws.Range("P17").Value = 0.11614513789011
$ws.Range("Q17").Value = 0.03008975795833333
$ws.Range("R17").Value = 0.270807821625
$ws.Range("S17").Value = 0.01606526076227096
$ws.Range("T17").Value = 0.01606526076227095

# Row 18
$ws.Range("G18").Value = 0.010746
$ws.Range("H18").Value = 0.032238
$ws.Range("I18").Value = 0.03164220884957143
$ws.Range("J18").Value = 0.03164220884957143
$ws.Range("M18").Value = 3.339352
$ws.Range("N18").Value = 10.018056
$ws.Range("O18").Value = 0.6054960700393903
$ws.Range("P18").Value = 0.6054960700393903
$ws.Range("Q18").Value = 0.035884676592
$ws.Range("R18").Value = 0.3229620893280001
$ws.Range("S18").Value = 0.01915923310578111
$ws.Range("T18").Value = 0.01915923310578111

# Row 19
$ws.Range("G19").Value = 0.010746
$ws.Range("H19").Value = 0.032238
$ws.Range("I19").Value = 0.03164220884957143
$ws.Range("J19").Value = 0.03164220884957143
$ws.Range("O19").Value = 0.2540955070726236
$ws.Range("P19").Value = 0.2540955070726236
$ws.Range("Q19").Value = 0.015058950084
$ws.Range("R19").Value = 0.135530550756
$ws.Range("S19").Value = 0.00804014310252971
$ws.Range("T19").Value = 0.008040143102529708

# Row 20
$ws.Range("G20").Value = 0.010746
$ws.Range("H20").Value = 0.032238
$ws.Range("I20").Value = 0.03164220884957143
$ws.Range("J20").Value = 0.03164220884957143
$ws.Range("K20").Value = 2
$ws.Range("L20").Value = 0.6666666666666666
$ws.Range("M20").Value = 0.1338136666666667
$ws.Range("N20").Value = 0.401441
$ws.Range("O20").Value = 0.02426328499787613
$ws.Range("P20").Value = 0.02426328499787612
$ws.Range("Q20").Value = 0.001437961662
$ws.Range("R20").Value = 0.012941654958
$ws.Range("S20").Value = 0.0007677439312794695
$ws.Range("T20").Value = 0.0007677439312794694

# Row 21
$ws.Range("G21").Value = 0.010746
$ws.Range("H21").Value = 0.032238
$ws.Range("I21").Value = 0.03164220884957143
$ws.Range("J21").Value = 0.03164220884957143
$ws.Range("M21").Value = 0.6405483333333334
$ws.Range("N21").Value = 1.921645
$ws.Range("O21").Value = 0.11614513789011
$ws.Range("P21").Value = 0.11614513789011
$ws.Range("Q21").Value = 0.006883332390000001
$ws.Range("R21").Value = 0.06194999151
$ws.Range("S21").Value = 0.003675088709981133
$ws.Range("T21").Value = 0.003675088709981133

# Row 22
$ws.Range("E22").Value = 3
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 0.03299933333333333
$ws.Range("H22").Value = 0.098998
$ws.Range("I22").Value = 0.09716841589707399
$ws.Range("J22").Value = 0.09716841589707399
$ws.Range("M22").Value = 3.339352
$ws.Range("N22").Value = 10.018056
$ws.Range("O22").Value = 0.6054960700393903
$ws.Range("P22").Value = 0.6054960700393903
$ws.Range("Q22").Value = 0.1101963897653333
$ws.Range("R22").Value = 0.9917675078880002
$ws.Range("S22").Value = 0.05883509395763132
$ws.Range("T22").Value = 0.05883509395763132

# Row 23
$ws.Range("E23").Value = 3
$ws.Range("F23").Value = 1
$ws.Range("G23").Value = 0.03299933333333333
$ws.Range("H23").Value = 0.098998
$ws.Range("I23").Value = 0.09716841589707399
$ws.Range("J23").Value = 0.09716841589707399
$ws.Range("O23").Value = 0.2540955070726236
$ws.Range("P23").Value = 0.2540955070726236
$ws.Range("Q23").Value = 0.046243747764
$ws.Range("R23").Value = 0.4161937298760001
$ws.Range("S23").Value = 0.0246900579088106
$ws.Range("T23").Value = 0.02469005790881059

# Row 24
$ws.Range("E24").Value = 3
$ws.Range("F24").Value = 1
$ws.Range("G24").Value = 0.03299933333333333
$ws.Range("H24").Value = 0.098998
$ws.Range("I24").Value = 0.09716841589707399
$ws.Range("J24").Value = 0.09716841589707399
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.1338136666666667
$ws.Range("N24").Value = 0.401441
$ws.Range("O24").Value = 0.02426328499787613
$ws.Range("P24").Value = 0.02426328499787612
$ws.Range("Q24").Value = 0.00441576179088889
$ws.Range("R24").Value = 0.039741856118
$ws.Range("S24").Value = 0.002357624967702864
$ws.Range("T24").Value = 0.002357624967702863

# Row 25
$ws.Range("E25").Value = 3
$ws.Range("F25").Value = 1
$ws.Range("G25").Value = 0.03299933333333333
$ws.Range("H25").Value = 0.098998
$ws.Range("I25").Value = 0.09716841589707399
$ws.Range("J25").Value = 0.09716841589707399
$ws.Range("M25").Value = 0.6405483333333334
$ws.Range("N25").Value = 1.921645
$ws.Range("O25").Value = 0.11614513789011
$ws.Range("P25").Value = 0.11614513789011
$ws.Range("Q25").Value = 0.02113766796777778
$ws.Range("R25").Value = 0.19023901171
$ws.Range("S25").Value = 0.01128563906292922
$ws.Range("T25").Value = 0.01128563906292922

